$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "50-1=",
    "52-30=",
    "0+71=",
    "32+44=",
    "94-79=",
    "15+40=",
    "26+59=",
    "9-8=",
    "99-49=",
    "49+40=",
    "17+41=",
    "11+69=",
    "70-36=",
    "68-56=",
    "55-52=",
    "68-29=",
    "75-25=",
    "65+9=",
    "51-14=",
    "16-3=",
    "96-14=",
    "56+35=",
    "67-15=",
    "51-9=",
    "7+26=",
    "63-50=",
    "61-47=",
    "88-65=",
    "38+42=",
    "44+27=",
    "87-76=",
    "40+57=",
    "42-7=",
    "8+82=",
    "43+34=",
    "85+9=",
    "43+55=",
    "63-5=",
    "57+28=",
    "87-80=",
    "87-49=",
    "46-15=",
    "52+16=",
    "99-74=",
    "40+5=",
    "62-14=",
    "66-12=",
    "98-97=",
    "67-64=",
    "68-63=",
    "91-2=",
    "35-21=",
    "5+54=",
    "15-9=",
    "3+27=",
    "72+6=",
    "12+69=",
    "72-22=",
    "90-52=",
    "30-12=",
    "61-56=",
    "87-6=",
    "41-21=",
    "19-13=",
    "70-52=",
    "69-32=",
    "58+38=",
    "80+0=",
    "15+64=",
    "99-93=",
    "30+10=",
    "13+39=",
    "28-21=",
    "13+34=",
    "85-28=",
    "75-56=",
    "74-13=",
    "20-5=",
    "95-41=",
    "54-44=",
    "44+45=",
    "37-13=",
    "34-18=",
    "38+21=",
    "92-88=",
    "19+24=",
    "49+17=",
    "78-64=",
    "69+29=",
    "21-21=",
    "58-5=",
    "34+17=",
    "99-13=",
    "73+23=",
    "31+44=",
    "34+13=",
    "48-43=",
    "49+21=",
    "58+6=",
    "79-43="
)

$cols = 5
$idx = 0
for ($row = 1; $row -le $t.Rows.Count; $row++) {
    for ($col = 1; $col -le $cols; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated" $idx "cells (expected" $newValues.Count ")"